# Control Arm Naming Guide - "Clevis and shim updates"
#
# Changes applied (per commit diff):
#   - Parts!D5 (Inboard Clevis, YY):        17 -> 22
#   - Parts!I5 (Inboard Clevis, Status):     "Carry Over - FEA OK" -> "Done"
#   - Parts!D6 (Inboard Clevis Shim, YY):    17 -> 22
#   - Parts!I6 (Inboard Clevis Shim, Status):"Carry Over" -> "Done"
#   - New shim-thickness/qty annotations added on row 6: L6:O6
#   - Minor view-state nudges (active cell selection) to mirror the saved file

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parts")

# Inboard Clevis (row 5): YY goes from 17 to 22, status flips to Done
$ws.Range("D5").Value = 22
$ws.Range("I5").Value = "Done"

# Inboard Clevis Shim (row 6): YY goes from 17 to 22, status flips to Done
$ws.Range("D6").Value = 22
$ws.Range("I6").Value = "Done"

# New shim breakdown notes added alongside the Inboard Clevis Shim row
$ws.Range("L6").Value = "0.3125 x 4"
$ws.Range("M6").Value = "0.125 x 20"
$ws.Range("N6").Value = "0.063 x 8"
$ws.Range("O6").Value = "0.032 x 8"

# Give the new columns a sensible width (matches the author's saved widths)
$ws.Columns.Item(12).ColumnWidth = 8.61
$ws.Columns.Item(13).ColumnWidth = 8.5

# Reflect the author's final cursor position in the saved view state
$ws.Range("H17").Select() | Out-Null
